$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I holds "% of Q Drop's" - add header + per-course values.
# Values must stay text (matching the existing "% of X's" columns), so
# force Text number format on each target cell before assigning the
# string value (otherwise "0.00%" etc. get auto-converted to numbers).

$cells = @{
    "I1"  = "% of Q Drop's"
    "I3"  = "0.00%"
    "I6"  = "0.00%"
    "I9"  = "0.00%"
    "I12" = "0.00%"
    "I15" = "0.00%"
    "I18" = "0.00%"
    "I21" = "15.38%"
    "I24" = "0.00%"
    "I27" = "0.00%"
    "I30" = "0.00%"
    "I33" = "0.00%"
    "I36" = "0.00%"
    "I39" = "0.00%"
    "I42" = "0.00%"
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}
